# Auto-generated edit script: update Leve market-price derived columns (H-N)
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 204.94118
$ws.Range("I33").Value = 128
$ws.Range("J33").Value = 455
$ws.Range("K33").Value = 128
$ws.Range("L33").Value = 455
$ws.Range("M33").Value = 101
$ws.Range("N33").Value = -913
$ws.Range("H53").Value = 478.66666
$ws.Range("I53").Value = 349.23077
$ws.Range("J53").Value = 631.63635
$ws.Range("K53").Value = 349.23077
$ws.Range("L53").Value = 631.63635
$ws.Range("M53").Value = 287.76923
$ws.Range("N53").Value = -1905.63635
$ws.Range("H112").Value = 19609434
$ws.Range("J112").Value = 1643.3125
$ws.Range("L112").Value = 4929.9375
$ws.Range("N112").Value = -7145.9375
$ws.Range("H129").Value = 760.2143
$ws.Range("I129").Value = 357.84616
$ws.Range("J129").Value = 1108.9333
$ws.Range("K129").Value = 1073.53848
$ws.Range("L129").Value = 3326.7999
$ws.Range("M129").Value = 3926.46152
$ws.Range("N129").Value = -13326.7999
$ws.Range("H132").Value = 23261220
$ws.Range("I132").Value = 30308248
$ws.Range("J132").Value = 6031
$ws.Range("K132").Value = 90924744
$ws.Range("L132").Value = 18093
$ws.Range("M132").Value = -90922214
$ws.Range("N132").Value = -23153
$ws.Range("H138").Value = 4662.2603
$ws.Range("J138").Value = 5326.427
$ws.Range("L138").Value = 15979.281
$ws.Range("N138").Value = -26259.281

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1254.931
$ws.Range("I61").Value = 907.4583
$ws.Range("J61").Value = 2922.8
$ws.Range("K61").Value = 907.4583
$ws.Range("L61").Value = 2922.8
$ws.Range("M61").Value = -695.4583
$ws.Range("N61").Value = -3346.8
$ws.Range("H132").Value = 2557.5
$ws.Range("I132").Value = 1623.9333
$ws.Range("J132").Value = 4891.4165
$ws.Range("K132").Value = 4871.7999
$ws.Range("L132").Value = 14674.2495
$ws.Range("M132").Value = -2341.7999
$ws.Range("N132").Value = -19734.2495
$ws.Range("H136").Value = 1254.931
$ws.Range("I136").Value = 907.4583
$ws.Range("J136").Value = 2922.8
$ws.Range("K136").Value = 2722.3749
$ws.Range("L136").Value = 8768.400000000001
$ws.Range("M136").Value = -172.3748999999998
$ws.Range("N136").Value = -13868.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1458.1428
$ws.Range("I86").Value = 1240
$ws.Range("J86").Value = 2003.5
$ws.Range("K86").Value = 1240
$ws.Range("L86").Value = 2003.5
$ws.Range("M86").Value = -117
$ws.Range("N86").Value = -4249.5
$ws.Range("H89").Value = 1458.1428
$ws.Range("I89").Value = 1240
$ws.Range("J89").Value = 2003.5
$ws.Range("K89").Value = 6200
$ws.Range("L89").Value = 10017.5
$ws.Range("M89").Value = -584
$ws.Range("N89").Value = -21249.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1680.6129
$ws.Range("I58").Value = 1579.5797
$ws.Range("J58").Value = 1971.0834
$ws.Range("K58").Value = 1579.5797
$ws.Range("L58").Value = 1971.0834
$ws.Range("M58").Value = -1376.5797
$ws.Range("N58").Value = -2377.0834
$ws.Range("H136").Value = 1680.6129
$ws.Range("I136").Value = 1579.5797
$ws.Range("J136").Value = 1971.0834
$ws.Range("K136").Value = 4738.7391
$ws.Range("L136").Value = 5913.2502
$ws.Range("M136").Value = -2188.7391
$ws.Range("N136").Value = -11013.2502

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 87.64706
$ws.Range("I12").Value = 21.5
$ws.Range("J12").Value = 123.72727
$ws.Range("K12").Value = 64.5
$ws.Range("L12").Value = 371.18181
$ws.Range("M12").Value = 108.5
$ws.Range("N12").Value = -717.18181
$ws.Range("H131").Value = 828.386
$ws.Range("I131").Value = 520.1111
$ws.Range("J131").Value = 886.1875
$ws.Range("K131").Value = 1560.3333
$ws.Range("L131").Value = 2658.5625
$ws.Range("M131").Value = 3479.6667
$ws.Range("N131").Value = -12738.5625
$ws.Range("H132").Value = 1967.7
$ws.Range("I132").Value = 994.5
$ws.Range("J132").Value = 2616.5
$ws.Range("K132").Value = 8950.5
$ws.Range("L132").Value = 23548.5
$ws.Range("M132").Value = -6420.5
$ws.Range("N132").Value = -28608.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 15627492
$ws.Range("I80").Value = 22729592
$ws.Range("J80").Value = 2869.2
$ws.Range("K80").Value = 22729592
$ws.Range("L80").Value = 2869.2
$ws.Range("M80").Value = -22728594
$ws.Range("N80").Value = -4865.2
$ws.Range("H83").Value = 15627492
$ws.Range("I83").Value = 22729592
$ws.Range("J83").Value = 2869.2
$ws.Range("K83").Value = 113647960
$ws.Range("L83").Value = 14346
$ws.Range("M83").Value = -113642968
$ws.Range("N83").Value = -24330
$ws.Range("H126").Value = 5171.283
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 5171.283
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 15513.849
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -20453.849
$ws.Range("H132").Value = 2312.4807
$ws.Range("I132").Value = 1495.4412
$ws.Range("J132").Value = 3855.7778
$ws.Range("K132").Value = 4486.3236
$ws.Range("L132").Value = 11567.3334
$ws.Range("M132").Value = -1956.3236
$ws.Range("N132").Value = -16627.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4358.685
$ws.Range("I132").Value = 1587.9286
$ws.Range("J132").Value = 7342.577
$ws.Range("K132").Value = 4763.7858
$ws.Range("L132").Value = 22027.731
$ws.Range("M132").Value = -2233.7858
$ws.Range("N132").Value = -27087.731
$ws.Range("H133").Value = 45366.6
$ws.Range("J133").Value = 45366.6
$ws.Range("L133").Value = 45366.6
$ws.Range("N133").Value = -50426.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 22961032
$ws.Range("I81").Value = 24727144
$ws.Range("J81").Value = 1600
$ws.Range("K81").Value = 49454288
$ws.Range("L81").Value = 3200
$ws.Range("M81").Value = -49453227
$ws.Range("N81").Value = -5322
$ws.Range("H84").Value = 22961032
$ws.Range("I84").Value = 24727144
$ws.Range("J84").Value = 1600
$ws.Range("K84").Value = 247271440
$ws.Range("L84").Value = 16000
$ws.Range("M84").Value = -247266136
$ws.Range("N84").Value = -26608
$ws.Range("H122").Value = 3074.3914
$ws.Range("I122").Value = 1777.1177
$ws.Range("J122").Value = 6750
$ws.Range("K122").Value = 5331.3531
$ws.Range("L122").Value = 20250
$ws.Range("M122").Value = -2881.3531
$ws.Range("N122").Value = -25150

